# First draft of presentation.
# 1) Append two new slides (Followup / Questions) after the existing 9 slides.
# 2) Re-format the "Possible code" slide's code block: switch the font of
#    every run to Consolas and drop one leading tab stop from every line
#    (the block now relies on the shape's own tab handling for the first
#    indent level instead of a literal tab).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 10: "Followup"
# ---------------------------------------------------------------------
$s10 = $p.Slides.Add(10, 2)
$s10.Shapes.Item(1).TextFrame.TextRange.Text = "Followup"

$body10 = $s10.Shapes.Item(2).TextFrame.TextRange
$body10.Text = "Can we make the code more efficient?`rIf the original data were collected by checked in, checked out pairs, how do we create the data we want?`rHow could we handle multiple books in the data?"

# ---------------------------------------------------------------------
# Slide 11: "Questions"
# ---------------------------------------------------------------------
$s11 = $p.Slides.Add(11, 2)
$s11.Shapes.Item(1).TextFrame.TextRange.Text = "Questions"

# ---------------------------------------------------------------------
# Slide 9: reformat the code sample
#   - switch every run in the code block to Consolas
#   - the block now opens one tab stop narrower, so drop a single
#     leading tab character from every line
# ---------------------------------------------------------------------
$codeSlide = $p.Slides.Item(9)
$codeBody = $codeSlide.Shapes.Item(2).TextFrame.TextRange

# number of leading tab characters currently on each of the 18 code lines
$leadingTabs = @(1, 2, 2, 2, 3, 3, 3, 4, 4, 3, 4, 3, 3, 4, 3, 2, 2, 1)

for ($i = 1; $i -le $leadingTabs.Length; $i++) {
    $para = $codeBody.Paragraphs($i, 1)

    $n = $leadingTabs[$i - 1]
    $oldTabs = "`t" * $n
    $newTabs = "`t" * ($n - 1)
    $para.Replace($oldTabs, $newTabs) | Out-Null

    $para.Font.Name = "Consolas"
}
